$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 16.21836188366955
$ws.Range("C2").Value = 8.605727679572531
$ws.Range("D2").Value = 8.811127173145866
$ws.Range("F2").Value = 34.09013796710386
$ws.Range("G2").Value = 3.662729381142114
$ws.Range("I2").Value = 25.19782175660043
$ws.Range("J2").Value = 10.31796117393695
$ws.Range("L2").Value = 11.8432712171455
$ws.Range("M2").Value = 16.6034853301032
$ws.Range("N2").Value = 18.86876288371683
$ws.Range("O2").Value = 25.86627466905471
$ws.Range("B3").Value = 15.83908569977678
$ws.Range("C3").Value = 8.355152980597151
$ws.Range("D3").Value = 8.818453654927906
$ws.Range("F3").Value = 34.15452684217377
$ws.Range("G3").Value = 3.664735672705961
$ws.Range("I3").Value = 25.29566029792337
$ws.Range("J3").Value = 10.33500864169642
$ws.Range("L3").Value = 11.84894389026018
$ws.Range("M3").Value = 16.52864988538082
$ws.Range("N3").Value = 18.92061666096053
$ws.Range("O3").Value = 25.92864198422926
$ws.Range("B4").Value = 15.60373641578065
$ws.Range("C4").Value = 8.196057950913618
$ws.Range("D4").Value = 8.82379069391115
$ws.Range("F4").Value = 34.20169187913827
$ws.Range("G4").Value = 3.666034153018252
$ws.Range("I4").Value = 25.36037291759616
$ws.Range("J4").Value = 10.34602133646689
$ws.Range("L4").Value = 11.8537564335078
$ws.Range("M4").Value = 16.48463637088586
$ws.Range("N4").Value = 18.95421606230185
$ws.Range("O4").Value = 25.97245907011188
$ws.Range("B5").Value = 15.5073441134522
$ws.Range("C5").Value = 8.129968578666034
$ws.Range("D5").Value = 8.826176939829532
$ws.Range("F5").Value = 34.22282753450633
$ws.Range("G5").Value = 3.666580095822614
$ws.Range("I5").Value = 25.3879096328917
$ws.Range("J5").Value = 10.35064665646128
$ws.Range("L5").Value = 11.85605274897534
$ws.Range("M5").Value = 16.46719978380673
$ws.Range("N5").Value = 18.96835191582451
$ws.Range("O5").Value = 25.99170124637415
$ws.Range("B6").Value = 15.49131300538517
$ws.Range("C6").Value = 8.118920466214275
$ws.Range("D6").Value = 8.826585953958038
$ws.Range("F6").Value = 34.22645268619743
$ws.Range("G6").Value = 3.666671765514148
$ws.Range("I6").Value = 25.39255247150075
$ws.Range("J6").Value = 10.35142300699014
$ws.Range("L6").Value = 11.85645431729993
$ws.Range("M6").Value = 16.4643349617114
$ws.Range("N6").Value = 18.97072599948045
$ws.Range("O6").Value = 25.9949800494951
$ws.Range("B7").Value = 15.60243821754633
$ws.Range("C7").Value = 8.195171650913402
$ws.Range("D7").Value = 8.82382201921857
$ws.Range("F7").Value = 34.20196917059702
$ws.Range("G7").Value = 3.666041447695945
$ws.Range("I7").Value = 25.36073956816507
$ws.Range("J7").Value = 10.34608315766305
$ws.Range("L7").Value = 11.85378604426664
$ws.Range("M7").Value = 16.48439917730589
$ws.Range("N7").Value = 18.95440490498681
$ws.Range("O7").Value = 25.97271296703665
$ws.Range("B8").Value = 16.08818039889321
$ws.Range("C8").Value = 8.520452332696662
$ws.Range("D8").Value = 8.813479607255276
$ws.Range("F8").Value = 34.11075391669249
$ws.Range("G8").Value = 3.6634073571495
$ws.Range("I8").Value = 25.2305931015487
$ws.Range("J8").Value = 10.323726162859
$ws.Range("L8").Value = 11.84495174364568
$ws.Range("M8").Value = 16.57728797458686
$ws.Range("N8").Value = 18.88627718554238
$ws.Range("O8").Value = 25.88663110847734
$ws.Range("B9").Value = 17.0154346532665
$ws.Range("C9").Value = 9.114324655561195
$ws.Range("D9").Value = 8.799829531639194
$ws.Range("F9").Value = 33.9925422191072
$ws.Range("G9").Value = 3.658768052300087
$ws.Range("I9").Value = 25.01222352772065
$ws.Range("J9").Value = 10.28419431800459
$ws.Range("L9").Value = 11.83813795108567
$ws.Range("M9").Value = 16.77422944746113
$ws.Range("N9").Value = 18.76660537101606
$ws.Range("O9").Value = 25.76175429489508
$ws.Range("B10").Value = 17.67416702692249
$ws.Range("C10").Value = 9.520943335576156
$ws.Range("D10").Value = 8.793813359415147
$ws.Range("F10").Value = 33.94280957716857
$ws.Range("G10").Value = 3.655676982386326
$ws.Range("I10").Value = 24.87429964776256
$ws.Range("J10").Value = 10.25775252421358
$ws.Range("L10").Value = 11.83948518092682
$ws.Range("M10").Value = 16.92713947082462
$ws.Range("N10").Value = 18.68710815868198
$ws.Range("O10").Value = 25.69691831817183
$ws.Range("B11").Value = 17.96759700058733
$ws.Range("C11").Value = 9.698964623846924
$ws.Range("D11").Value = 8.791940935479433
$ws.Range("F11").Value = 33.92826260223183
$ws.Range("G11").Value = 3.654338995213203
$ws.Range("I11").Value = 24.81645377570235
$ws.Range("J11").Value = 10.24628313123315
$ws.Range("L11").Value = 11.84146490831906
$ws.Range("M11").Value = 16.99830757357049
$ws.Range("N11").Value = 18.65275893588424
$ws.Range("O11").Value = 25.6732888376757
$ws.Range("B12").Value = 18.07771309171482
$ws.Range("C12").Value = 9.765340222536862
$ws.Range("D12").Value = 8.791355573886568
$ws.Range("F12").Value = 33.92391585926809
$ws.Range("G12").Value = 3.653842081537641
$ws.Range("I12").Value = 24.79525417435037
$ws.Range("J12").Value = 10.2420199665705
$ws.Range("L12").Value = 11.84240992399266
$ws.Range("M12").Value = 17.02547267635486
$ws.Range("N12").Value = 18.64001171884683
$ws.Range("O12").Value = 25.66518566088816
$ws.Range("B13").Value = 18.05404363613635
$ws.Range("C13").Value = 9.751091733886733
$ws.Range("D13").Value = 8.791476151031432
$ws.Range("F13").Value = 33.92480033242682
$ws.Range("G13").Value = 3.653948667775898
$ws.Range("I13").Value = 24.79978849557287
$ws.Range("J13").Value = 10.2429345609159
$ws.Range("L13").Value = 11.84219773046829
$ws.Range("M13").Value = 17.01961286882286
$ws.Range("N13").Value = 18.64274550561065
$ws.Range("O13").Value = 25.66689323268372
$ws.Range("B14").Value = 17.97667691303583
$ws.Range("C14").Value = 9.704446371844245
$ws.Range("D14").Value = 8.791890302179057
$ws.Range("F14").Value = 33.92788170680695
$ws.Range("G14").Value = 3.654297918633202
$ws.Range("I14").Value = 24.81469552760534
$ws.Range("J14").Value = 10.24593079589536
$ws.Range("L14").Value = 11.84153874772525
$ws.Range("M14").Value = 17.00053824111396
$ws.Range("N14").Value = 18.6517050067906
$ws.Range("O14").Value = 25.67260524627885
$ws.Range("B15").Value = 17.92915446506435
$ws.Range("C15").Value = 9.675738622375423
$ws.Range("D15").Value = 8.792160070618078
$ws.Range("F15").Value = 33.92992045191556
$ws.Range("G15").Value = 3.654513113589869
$ws.Range("I15").Value = 24.82391841397175
$ws.Range("J15").Value = 10.24777649017394
$ws.Range("L15").Value = 11.84116050294711
$ws.Range("M15").Value = 16.98888203680728
$ws.Range("N15").Value = 18.65722680685585
$ws.Range("O15").Value = 25.67621407458575
$ws.Range("B16").Value = 17.65485628001236
$ws.Range("C16").Value = 9.509166035093063
$ws.Range("D16").Value = 8.793953068347173
$ws.Range("F16").Value = 33.9439228334036
$ws.Range("G16").Value = 3.655765790379753
$ws.Range("I16").Value = 24.87817864860665
$ws.Range("J16").Value = 10.25851329643782
$ws.Range("L16").Value = 11.83938320339769
$ws.Range("M16").Value = 16.92251947388481
$ws.Range("N16").Value = 18.68938940142417
$ws.Range("O16").Value = 25.69858070496325
$ws.Range("B17").Value = 17.48491121390845
$ws.Range("C17").Value = 9.405170613455754
$ws.Range("D17").Value = 8.795273972254769
$ws.Range("F17").Value = 33.95458199690852
$ws.Range("G17").Value = 3.65655168932699
$ws.Range("I17").Value = 24.91272053549409
$ws.Range("J17").Value = 10.2652429283086
$ws.Range("L17").Value = 11.83864213026033
$ws.Range("M17").Value = 16.88220888487088
$ws.Range("N17").Value = 18.7095842418717
$ws.Range("O17").Value = 25.71380513825445
$ws.Range("B18").Value = 17.38658252495902
$ws.Range("C18").Value = 9.344702851420521
$ws.Range("D18").Value = 8.796115117165817
$ws.Range("F18").Value = 33.96147308998375
$ws.Range("G18").Value = 3.657010135611748
$ws.Range("I18").Value = 24.93304890954261
$ws.Range("J18").Value = 10.26916628176799
$ws.Range("L18").Value = 11.83834461768414
$ws.Range("M18").Value = 16.85917568652686
$ws.Range("N18").Value = 18.72137061361437
$ws.Range("O18").Value = 25.72311377808166
$ws.Range("B19").Value = 17.35319365554735
$ws.Range("C19").Value = 9.324118678503563
$ws.Range("D19").Value = 8.796413912386992
$ws.Range("F19").Value = 33.96393684044786
$ws.Range("G19").Value = 3.657166461404554
$ws.Range("I19").Value = 24.94001085079442
$ws.Range("J19").Value = 10.2705037148835
$ws.Range("L19").Value = 11.83826603104482
$ws.Range("M19").Value = 16.85140369684059
$ws.Range("N19").Value = 18.72539064298925
$ws.Range("O19").Value = 25.72636027195657
$ws.Range("B20").Value = 17.50306304019333
$ws.Range("C20").Value = 9.416308904336194
$ws.Range("D20").Value = 8.79512493995372
$ws.Range("F20").Value = 33.95336862501588
$ws.Range("G20").Value = 3.656467365150555
$ws.Range("I20").Value = 24.90899579293884
$ws.Range("J20").Value = 10.26452110129522
$ws.Range("L20").Value = 11.83870770397336
$ws.Range("M20").Value = 16.8864843595486
$ws.Range("N20").Value = 18.70741679204641
$ws.Range("O20").Value = 25.71212733522017
$ws.Range("B21").Value = 17.99942930564064
$ws.Range("C21").Value = 9.71817567361675
$ws.Range("D21").Value = 8.791765304157167
$ws.Range("F21").Value = 33.92694509873835
$ws.Range("G21").Value = 3.654195070874047
$ws.Range("I21").Value = 24.81029781701717
$ws.Range("J21").Value = 10.24504855944496
$ws.Range("L21").Value = 11.84172701561834
$ws.Range("M21").Value = 17.00613520802561
$ws.Range("N21").Value = 18.64906633190249
$ws.Range("O21").Value = 25.670904551895
$ws.Range("B22").Value = 18.31796527594694
$ws.Range("C22").Value = 9.909402890602927
$ws.Range("D22").Value = 8.790290217426543
$ws.Range("F22").Value = 33.9164478973756
$ws.Range("G22").Value = 3.652766821454446
$ws.Range("I22").Value = 24.74990505082932
$ws.Range("J22").Value = 10.23278853230775
$ws.Range("L22").Value = 11.84483826368206
$ws.Range("M22").Value = 17.08558117527733
$ws.Range("N22").Value = 18.61244655511492
$ws.Range("O22").Value = 25.64888764300084
$ws.Range("B23").Value = 18.14852569770217
$ws.Range("C23").Value = 9.807907189790988
$ws.Range("D23").Value = 8.791011768891789
$ws.Range("F23").Value = 33.92143081219032
$ws.Range("G23").Value = 3.653523921109988
$ws.Range("I23").Value = 24.78176112674346
$ws.Range("J23").Value = 10.23928937961974
$ws.Range("L23").Value = 11.8430740293518
$ws.Range("M23").Value = 17.04307061591788
$ws.Range("N23").Value = 18.63185280843277
$ws.Range("O23").Value = 25.66018748422173
$ws.Range("B24").Value = 17.49485853831654
$ws.Range("C24").Value = 9.411275396383497
$ws.Range("D24").Value = 8.795192062813301
$ws.Range("F24").Value = 33.95391481379582
$ws.Range("G24").Value = 3.65650546747843
$ws.Range("I24").Value = 24.91067828578396
$ws.Range("J24").Value = 10.2648472698148
$ws.Range("L24").Value = 11.83867765761953
$ws.Range("M24").Value = 16.88455097423082
$ws.Range("N24").Value = 18.70839614753413
$ws.Range("O24").Value = 25.71288413844646
$ws.Range("B25").Value = 16.76805713004594
$ws.Range("C25").Value = 8.958693064026111
$ws.Range("D25").Value = 8.802815328760305
$ws.Range("F25").Value = 34.01801107489409
$ws.Range("G25").Value = 3.6599671251641
$ws.Range("I25").Value = 25.06734851585044
$ws.Range("J25").Value = 10.29442992816667
$ws.Range("L25").Value = 11.83813795108567
$ws.Range("M25").Value = 16.77422944746113
$ws.Range("N25").Value = 18.76660537101606
$ws.Range("O25").Value = 25.76175429489508
